$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '45.945.39'
Set-TextValue 'E2' '  +6.94%  '
Set-TextValue 'D3' '2.424.56'
Set-TextValue 'E3' '  +5.82%  '
Set-TextValue 'E4' '  -0.25%  '
Set-TextValue 'D5' '116.72'
Set-TextValue 'E5' '  +13.36%  '
Set-TextValue 'D6' '320.09'
Set-TextValue 'E6' '  +2.73%  '
Set-TextValue 'E7' '  +3.20%  '
Set-TextValue 'E8' '  +0.01%  '
Set-TextValue 'D9' '0.634'
Set-TextValue 'E9' '  +5.05%  '
Set-TextValue 'D10' '43.25'
Set-TextValue 'E10' '  +10.92%  '
Set-TextValue 'E11' '  +4.53%  '
Set-TextValue 'D12' '8.76'
Set-TextValue 'E12' '  +6.35%  '
Set-TextValue 'E14' '  +2.30%  '
Set-TextValue 'E15' '  +4.59%  '
Set-TextValue 'D16' '2.792.22'
Set-TextValue 'D17' '2.429.80'
Set-TextValue 'E17' '  +5.87%  '
Set-TextValue 'D18' '45.797.53'
Set-TextValue 'E18' '  +7.09%  '
Set-TextValue 'E19' '  +4.96%  '
Set-TextValue 'E20' '  +4.79%  '
Set-TextValue 'E21' '  -0.97%  '
Set-TextValue 'D22' '75.29'
Set-TextValue 'E22' '  +2.54%  '
Set-TextValue 'E23' '  +4.70%  '
Set-TextValue 'D24' '269.83'
Set-TextValue 'E24' '  +0.65%  '
Set-TextValue 'E25' '  +8.72%  '
Set-TextValue 'E26' '  +0.06%  '
Set-TextValue 'D27' '7.69'
Set-TextValue 'E27' '  +6.61%  '
Set-TextValue 'D28' '11.40'
Set-TextValue 'E28' '  +5.77%  '
Set-TextValue 'E29' '  +2.53%  '
Set-TextValue 'D30' '40.05'
Set-TextValue 'E30' '  +11.63%  '
Set-TextValue 'D31' '23.14'
Set-TextValue 'E31' '  +3.54%  '
Set-TextValue 'E32' '  +13.27%  '
Set-TextValue 'D33' '174.09'
Set-TextValue 'E33' '  +5.71%  '
Set-TextValue 'E34' '  +15.59%  '
Set-TextValue 'E35' '  +9.14%  '
Set-TextValue 'E36' '  +2.09%  '
Set-TextValue 'E37' '  +10.49%  '
Set-TextValue 'D38' '3.16'
Set-TextValue 'E38' '  +12.74%  '
Set-TextValue 'D39' '4.18'
Set-TextValue 'E39' '  +15.14%  '
Set-TextValue 'D40' '0.0368'
Set-TextValue 'E40' '  +6.14%  '
Set-TextValue 'E41' '  +17.30%  '
Set-TextValue 'D42' '102.38'
Set-TextValue 'E42' '  -3.90%  '
Set-TextValue 'B43' 'Celestia'
Set-TextValue 'C43' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D43' '13.72'
Set-TextValue 'E43' '  +13.70%  '
Set-TextValue 'B44' 'Algorand'
Set-TextValue 'C44' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.241'
Set-TextValue 'E44' '  +6.44%  '
Set-TextValue 'D45' '72.50'
Set-TextValue 'E45' '  +2.56%  '
Set-TextValue 'E46' '  -0.64%  '
Set-TextValue 'E47' '  +14.52%  '
Set-TextValue 'D48' '117.87'
Set-TextValue 'E48' '  +6.88%  '
Set-TextValue 'E49' '  +16.93%  '
Set-TextValue 'D51' '79.74'
Set-TextValue 'E51' '  +3.41%  '
